$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: make a cell a "text placeholder" cell (used for the " 0" / "***.*"
# style blank markers in this report) by copying both the number-format and
# the text value from an existing placeholder cell of the same kind. This
# avoids creating brand-new cell styles and reuses the existing shared
# string entries ("0" / "***.*").
# ---------------------------------------------------------------------------
function Set-Placeholder {
    param([string]$targetCell, [string]$sourceCell)
    $ws.Range($sourceCell).Copy()
    $ws.Range($targetCell).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($sourceCell).Copy()
    $ws.Range($targetCell).PasteSpecial(-4163)  # xlPasteValues
}

# ---------------------------------------------------------------------------
# Helper: turn a "text placeholder" cell back into a real numeric cell,
# reusing the number format of a stable, untouched cell of the desired kind
# (count style "C16" / percentage style "K17") so no new cell style gets
# created.
# ---------------------------------------------------------------------------
function Set-NumberFromPlaceholder {
    param([string]$targetCell, [string]$formatSourceCell, [double]$value)
    $ws.Range($formatSourceCell).Copy()
    $ws.Range($targetCell).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($targetCell).Value = $value
}

# ---------------------------------------------------------------------------
# Header text updates
# ---------------------------------------------------------------------------
# Mayor name
$ws.Range("M6").Value = "Thomas G. Donlon"

# "Volume 31   Number  38" -> "...39" (only the last run, chars 21-22, changes)
$ws.Range("A8").Characters(21, 2).Text = "39"

# "Report Covering the Week  9/16/2024  Through  9/22/2024"
#  -> "...9/23/2024  Through  9/29/2024"
$ws.Range("C9").Characters(27, 9).Text = "9/23/2024"
$ws.Range("C9").Characters(47, 9).Text = "9/29/2024"

# ---------------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------------
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 36
$ws.Range("K16").Value = 111.764705882353
$ws.Range("L16").Value = 71.428571428571
$ws.Range("M16").Value = 63.636363636363
$ws.Range("N16").Value = -78.048780487804

# ---------------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------------
$ws.Range("N17").Value = -65.625

# ---------------------------------------------------------------------------
# Row 18 (Burglary): D18 and E18 become blank placeholders ("0" / "***.*")
# ---------------------------------------------------------------------------
Set-Placeholder "D18" "C14"
Set-Placeholder "E18" "E14"
$ws.Range("N18").Value = -96.296296296296

# ---------------------------------------------------------------------------
# Row 19 (Gr. Larceny): C19 becomes a blank placeholder; D19/E19 become real
# numbers.
# ---------------------------------------------------------------------------
Set-Placeholder "C19" "C14"
Set-NumberFromPlaceholder "D19" "C16" 2
Set-NumberFromPlaceholder "E19" "K17" -100
$ws.Range("J19").Value = 40
$ws.Range("K19").Value = -10
$ws.Range("M19").Value = -37.931034482758
$ws.Range("N19").Value = -74.647887323943

# ---------------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = -50
$ws.Range("F21").Value = 4
$ws.Range("G21").Value = 7
$ws.Range("H21").Value = -42.857142857142
$ws.Range("I21").Value = 86
$ws.Range("J21").Value = 69
$ws.Range("K21").Value = 24.637681159420
$ws.Range("L21").Value = 50.877192982456
$ws.Range("M21").Value = -6.521739130434
$ws.Range("N21").Value = -77.486910994764

# ---------------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------------
$ws.Range("F24").Value = 7
$ws.Range("H24").Value = 40
$ws.Range("I24").Value = 32
$ws.Range("J24").Value = 32
$ws.Range("L24").Value = 33.333333333333
$ws.Range("M24").Value = -51.515151515151

# ---------------------------------------------------------------------------
# Row 26 (Misd. Assault): D26/E26 become blank placeholders.
# ---------------------------------------------------------------------------
Set-Placeholder "D26" "C14"
Set-Placeholder "E26" "E14"
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 31
$ws.Range("K26").Value = -24.390243902439
$ws.Range("L26").Value = 24
$ws.Range("M26").Value = 82.352941176470

# ---------------------------------------------------------------------------
# Row 27 (UCR Rape*): C27 and F27 become real numbers.
# ---------------------------------------------------------------------------
Set-NumberFromPlaceholder "C27" "C16" 1
Set-NumberFromPlaceholder "F27" "C16" 1
$ws.Range("I27").Value = 4
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = 100

# ---------------------------------------------------------------------------
# Row 28 (Other Sex Crimes): C28 becomes a blank placeholder.
# ---------------------------------------------------------------------------
Set-Placeholder "C28" "C14"

Write-Host "Edits applied"
